# Update the duty-roster statistics on sheet "历史" with the newly
# recalculated counts (values doubled / adjusted per the latest data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("历史")

# r => @(B, C, D, E, F, G, H, I)
$data = @{
    2  = @(4, 0, 4, 0,   2, 0, 2, 100)
    3  = @(4, 0, 4, 0,   0, 2, 2, 0)
    4  = @(6, 0, 6, 0,   0, 0, 0, 0)
    5  = @(0, 2, 2, 100, 0, 0, 0, 0)
    6  = @(0, 4, 4, 100, 2, 0, 2, 100)
    7  = @(0, 2, 2, 100, 0, 2, 2, 0)
    8  = @(2, 0, 2, 0,   2, 2, 4, 50)
    9  = @(2, 0, 2, 0,   2, 2, 4, 50)
    10 = @(2, 0, 2, 0,   0, 2, 2, 0)
    11 = @(2, 0, 2, 0,   2, 0, 2, 100)
    12 = @(4, 0, 4, 0,   0, 0, 0, 0)
    13 = @(2, 0, 2, 0,   0, 0, 0, 0)
    14 = @(2, 0, 2, 0,   0, 0, 0, 0)
    15 = @(2, 0, 2, 0,   0, 0, 0, 0)
    16 = @(2, 0, 2, 0,   0, 0, 0, 0)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
    $ws.Cells.Item($r, 8).Value = $vals[6]
    $ws.Cells.Item($r, 9).Value = $vals[7]
}
